$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 27-32: "Visualisation of court window" (VISU / RB) requirements ---
# Req ID | Type | Description | Status | Ticket ID | Component | Responsible | Planned For

$ws.Range("A27").Value = "UNIDEB_25"
$ws.Range("B27").Value = "H"
$ws.Range("C27").Value = "Visualisation of court window"
$ws.Range("D27").Value = "New"
$ws.Range("E27").Value = "https://trello.com/c/RpCLefE2/16-rd-visualisation"
$ws.Range("F27").Value = "VISU"
$ws.Range("G27").Value = "RB"
$ws.Range("H27").Value = "Sprint 1"

$ws.Range("A28").Value = "UNIDEB_26"
$ws.Range("B28").Value = "R"
$ws.Range("C28").Value = "The window displaying the court shall display the whole map (background bitmap)."
$ws.Range("D28").Value = "New"
$ws.Range("E28").Value = "https://trello.com/c/RpCLefE2/16-rd-visualisation"
$ws.Range("F28").Value = "VISU"
$ws.Range("G28").Value = "RB"
$ws.Range("H28").Value = "Sprint 1"

$ws.Range("A29").Value = "UNIDEB_27"
$ws.Range("B29").Value = "R"
$ws.Range("C29").Value = "The window should be resizable."
$ws.Range("D29").Value = "New"
$ws.Range("E29").Value = "https://trello.com/c/RpCLefE2/16-rd-visualisation"
$ws.Range("F29").Value = "VISU"
$ws.Range("G29").Value = "RB"
$ws.Range("H29").Value = "Sprint 1"

$ws.Range("A30").Value = "UNIDEB_28"
$ws.Range("B30").Value = "R"
$ws.Range("C30").Value = "The court (map) shall keep its ratio in the window, but fit to the window (in larger dimension of the map)."
$ws.Range("D30").Value = "New"
$ws.Range("E30").Value = "https://trello.com/c/RpCLefE2/16-rd-visualisation"
$ws.Range("F30").Value = "VISU"
$ws.Range("G30").Value = "RB"
$ws.Range("H30").Value = "Sprint 1"
$ws.Rows.Item(30).RowHeight = 30

$ws.Range("A31").Value = "UNIDEB_29"
$ws.Range("B31").Value = "R"
$ws.Range("C31").Value = "Blank area should be filled with solid (default) color."
$ws.Range("D31").Value = "New"
$ws.Range("E31").Value = "https://trello.com/c/RpCLefE2/16-rd-visualisation"
$ws.Range("F31").Value = "VISU"
$ws.Range("G31").Value = "RB"
$ws.Range("H31").Value = "Sprint 1"

$ws.Range("A32").Value = "UNIDEB_30"
$ws.Range("B32").Value = "R"
$ws.Range("C32").Value = "At start up, the window should be initialized with the window size of 800x600 px."
$ws.Range("D32").Value = "New"
$ws.Range("E32").Value = "https://trello.com/c/RpCLefE2/16-rd-visualisation"
$ws.Range("F32").Value = "VISU"
$ws.Range("G32").Value = "RB"
$ws.Range("H32").Value = "Sprint 1"

# --- View state: scroll the frozen sheet down and select C32 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$ws.Range("C32").Select()
